$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new row at position 3 (shifts current rows 3-8 down to 4-9),
# pushing "Ducks Unlimited" etc. down to make room for the new organization.
# (Cell content/styles shift with the insert; the worksheet's <hyperlinks>
# list does not, so that is repaired explicitly below.)
$ws.Rows(3).Insert()

# --- Row 2: Coalition for the Poudre River Watershed - updated coordinates ---
$ws.Range("G2").Value = -105.071971
$ws.Range("H2").Value = 40.596536999999998

# --- Row 3 (new): CSU Environmental Learning Center ---
$ws.Range("A3").Value = "CSU Environmental Learning Center"
$ws.Range("B3").Value = "University"
$ws.Range("C3").Value = "Environmental learning center"
$ws.Range("D3").Value = "Natural resources documentation"
$ws.Range("E3").Value = "https://warnercnr.colostate.edu/elc/"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = -105.019846
$ws.Range("H3").Value = 40.556621999999997

# --- Row 4: Ducks Unlimited - updated coordinates ---
$ws.Range("G4").Value = -105.028807
$ws.Range("H4").Value = 40.562764999999999

# --- Row 6: Rocky Mountain Flycasters and Trout Unlimited - updated coordinates ---
$ws.Range("G6").Value = -105.112968
$ws.Range("H6").Value = 40.402273000000001

# --- Repair the hyperlinks collection: the row insert shifted cell content
# down one row but left the old <hyperlink ref=.../> anchors in place, so
# clear them out and re-add one per surviving link, each still pointing at
# its original target URL but anchored on its new (shifted) cell. The new
# CSU row's E3 (plain text, no live link) intentionally gets none.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.poudrewatershed.org/")
$ws.Hyperlinks.Add($ws.Range("E7"), "http://www.savethepoudre.org/")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.fortcollinsdu.org/")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.rockymtnflycasters.org/")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.fcgov.com/naturalareas/")
$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.nature.org/en-us/get-involved/how-to-help/places-we-protect/phantom-canyon-preserve/")
$ws.Hyperlinks.Add($ws.Range("K9"), "https://www.nature.org/en-us/")
$ws.Hyperlinks.Add($ws.Range("K8"), "http://plattebasintimelapse.com/")
$ws.Hyperlinks.Add($ws.Range("E8"), "http://plattebasintimelapse.com/explore/galleries/north-fork-cache-la-poudre-river/")

# Re-apply the workbook's "Hyperlink" cell style (Hyperlinks.Add creates a
# fresh duplicate style; reset back to the shared built-in style already
# used throughout the sheet) on every linked cell, and also on the new
# CSU row's E3, which keeps the same visual style without a live link.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"
$ws.Range("E7").Style = "Hyperlink"
$ws.Range("E8").Style = "Hyperlink"
$ws.Range("E9").Style = "Hyperlink"
$ws.Range("K8").Style = "Hyperlink"
$ws.Range("K9").Style = "Hyperlink"

# --- Sheet view: scroll to column E, select H7 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H7").Select()
